$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows into the endpoint table: one for the new
# "DELETE /api/rooms/:id" entry (before current row 11) and one for the
# new "DELETE /api/tenants/:id" entry (before current row 16, i.e. right
# before the "/api/tenants/create" row once the first insert has shifted
# everything down by one).
$ws.Rows(11).Insert()
$ws.Rows(16).Insert()

# Header row
$ws.Range("B1").Value = "Method/Input inside body"

# "/"  (GET, serves html)
$ws.Range("B3").Value = "get"

# "/api"  (GET)
$ws.Range("B5").Value = "get"
$ws.Range("C5").Value = "{data:{message:”/api accessible”}}"

# "/api/login"  (POST)
$ws.Range("B7").Value = "post/{username,password}"

# "/api/rooms"  (GET)
$ws.Range("B9").Value = "get"

# "/api/rooms/:id"  (GET)
$ws.Range("B10").Value = "get"

# new row: "/api/rooms/:id"  (DELETE)
$ws.Range("A11").Value = "/api/rooms/:id"
$ws.Range("B11").Value = "delete"
$ws.Range("C11").Value = "{data:{message:”success”}}"
$ws.Range("D11").Value = "note the special delete method"

# "/api/rooms/create"  (POST)
$ws.Range("B12").Value = "post/{number,baseRent}"

# new row: "/api/tenants/:id"  (DELETE)
$ws.Range("A16").Value = "/api/tenants/:id"
$ws.Range("B16").Value = "delete"
$ws.Range("C16").Value = "{data:{message:”success”}}"

# "/api/tenants/create"  (POST, unchanged data already shifted into place)
$ws.Range("B17").Value = "{name,phoneNumber,aadharCard,room}"

# "/api/transactions/?room"  (GET)
$ws.Range("B19").Value = "get"

# "/api/transactions/:id"  (GET)
$ws.Range("B20").Value = "get"

# "/api/tenants/create/?roomNumber"  (POST)
$ws.Range("B21").Value = "post/{room,balance,transfer,remarks}"

# Restore the active selection recorded in the saved workbook
$ws.Range("D6").Select()
